$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 112
$ws.Range("H112").Value = 1457.6571
$ws.Range("I112").Value = 634.2857
$ws.Range("J112").Value = 1663.5
$ws.Range("K112").Value = 1902.8571
$ws.Range("L112").Value = 4990.5
$ws.Range("M112").Value = -794.8571000000002
$ws.Range("N112").Value = -7206.5

# Row 132
$ws.Range("H132").Value = 989.46155
$ws.Range("I132").Value = 783.32355
$ws.Range("J132").Value = 2391.2
$ws.Range("K132").Value = 2349.97065
$ws.Range("L132").Value = 7173.599999999999
$ws.Range("M132").Value = 180.0293500000002
$ws.Range("N132").Value = -12233.6

# Row 133
$ws.Range("H133").Value = 61000
$ws.Range("J133").Value = 61000
$ws.Range("L133").Value = 61000
$ws.Range("N133").Value = -71120

# Row 135
$ws.Range("H135").Value = 2143.5107
$ws.Range("I135").Value = 972.875
$ws.Range("J135").Value = 8832.857
$ws.Range("K135").Value = 8755.875
$ws.Range("L135").Value = 79495.713
$ws.Range("M135").Value = -6220.875
$ws.Range("N135").Value = -84565.713

# Row 141
$ws.Range("H141").Value = 1761.1489
$ws.Range("I141").Value = 1826.2858
$ws.Range("J141").Value = 1214
$ws.Range("K141").Value = 5478.857400000001
$ws.Range("L141").Value = 3642
$ws.Range("M141").Value = -298.8574000000008
$ws.Range("N141").Value = -14002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 306636.28
$ws.Range("I32").Value = 2152.057
$ws.Range("J32").Value = 4772404.5
$ws.Range("K32").Value = 2152.057
$ws.Range("L32").Value = 4772404.5
$ws.Range("M32").Value = -1865.057
$ws.Range("N32").Value = -4772978.5

# Row 61
$ws.Range("H61").Value = 963.4545000000001
$ws.Range("I61").Value = 628.45
$ws.Range("J61").Value = 1856.8
$ws.Range("K61").Value = 628.45
$ws.Range("L61").Value = 1856.8
$ws.Range("M61").Value = -416.45
$ws.Range("N61").Value = -2280.8

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 136
$ws.Range("H136").Value = 963.4545000000001
$ws.Range("I136").Value = 628.45
$ws.Range("J136").Value = 1856.8
$ws.Range("K136").Value = 1885.35
$ws.Range("L136").Value = 5570.4
$ws.Range("M136").Value = 664.6499999999999
$ws.Range("N136").Value = -10670.4

# Row 138
$ws.Range("H138").Value = 24950
$ws.Range("J138").Value = 24950
$ws.Range("L138").Value = 24950
$ws.Range("N138").Value = -35230

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 134
$ws.Range("H134").Value = 6623.5186
$ws.Range("I134").Value = 1353.3043
$ws.Range("J134").Value = 36927.25
$ws.Range("K134").Value = 4059.9129
$ws.Range("L134").Value = 110781.75
$ws.Range("M134").Value = -1524.9129
$ws.Range("N134").Value = -115851.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 6759413
$ws.Range("I31").Value = 7693701
$ws.Range("J31").Value = 11777.777
$ws.Range("K31").Value = 7693701
$ws.Range("L31").Value = 11777.777
$ws.Range("M31").Value = -7693406
$ws.Range("N31").Value = -12367.777

# Row 34
$ws.Range("H34").Value = 6759413
$ws.Range("I34").Value = 7693701
$ws.Range("J34").Value = 11777.777
$ws.Range("K34").Value = 7693701
$ws.Range("L34").Value = 11777.777
$ws.Range("M34").Value = -7693499
$ws.Range("N34").Value = -12181.777

# Row 58
$ws.Range("H58").Value = 1102.1945
$ws.Range("I58").Value = 917.63635
$ws.Range("J58").Value = 1392.2142
$ws.Range("K58").Value = 917.63635
$ws.Range("L58").Value = 1392.2142
$ws.Range("M58").Value = -714.63635
$ws.Range("N58").Value = -1798.2142

# Row 132
$ws.Range("H132").Value = 1286.2264
$ws.Range("I132").Value = 1127.0488
$ws.Range("J132").Value = 1830.0834
$ws.Range("K132").Value = 3381.1464
$ws.Range("L132").Value = 5490.2502
$ws.Range("M132").Value = -851.1464000000001
$ws.Range("N132").Value = -10550.2502

# Row 134
$ws.Range("H134").Value = 1334.4572
$ws.Range("I134").Value = 1359.5625
$ws.Range("J134").Value = 1066.6666
$ws.Range("K134").Value = 4078.6875
$ws.Range("L134").Value = 3199.9998
$ws.Range("M134").Value = -1543.6875
$ws.Range("N134").Value = -8269.9998

# Row 136
$ws.Range("H136").Value = 1102.1945
$ws.Range("I136").Value = 917.63635
$ws.Range("J136").Value = 1392.2142
$ws.Range("K136").Value = 2752.90905
$ws.Range("L136").Value = 4176.642599999999
$ws.Range("M136").Value = -202.9090500000002
$ws.Range("N136").Value = -9276.642599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 92
$ws.Range("H92").Value = 350.75
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 267.66666
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 802.9999799999999
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -3298.99998

# Row 122
$ws.Range("H122").Value = 626.3043
$ws.Range("I122").Value = 372.22223
$ws.Range("J122").Value = 789.6429000000001
$ws.Range("K122").Value = 3350.00007
$ws.Range("L122").Value = 7106.7861
$ws.Range("M122").Value = -900.0000700000001
$ws.Range("N122").Value = -12006.7861

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 132
$ws.Range("H132").Value = 1617.45
$ws.Range("I132").Value = 1624.1936
$ws.Range("J132").Value = 1594.2222
$ws.Range("K132").Value = 4872.5808
$ws.Range("L132").Value = 4782.6666
$ws.Range("M132").Value = -2342.5808
$ws.Range("N132").Value = -9842.6666

# Row 133
$ws.Range("H133").Value = 51897.777
$ws.Range("J133").Value = 51897.777
$ws.Range("L133").Value = 51897.777
$ws.Range("N133").Value = -62017.777

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 651.7213
$ws.Range("I22").Value = 528.5714
$ws.Range("J22").Value = 923.9474
$ws.Range("K22").Value = 528.5714
$ws.Range("L22").Value = 923.9474
$ws.Range("M22").Value = -233.5714
$ws.Range("N22").Value = -1513.9474

# Row 27
$ws.Range("H27").Value = 651.7213
$ws.Range("I27").Value = 528.5714
$ws.Range("J27").Value = 923.9474
$ws.Range("K27").Value = 528.5714
$ws.Range("L27").Value = 923.9474
$ws.Range("M27").Value = -421.5714
$ws.Range("N27").Value = -1137.9474

# Row 55
$ws.Range("H55").Value = 653.6667
$ws.Range("I55").Value = 656.13043
$ws.Range("J55").Value = 648
$ws.Range("K55").Value = 656.13043
$ws.Range("L55").Value = 648
$ws.Range("M55").Value = -483.13043
$ws.Range("N55").Value = -994

# Row 129
$ws.Range("H129").Value = 27000
$ws.Range("J129").Value = 27000
$ws.Range("L129").Value = 27000
$ws.Range("N129").Value = -37000

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 129
$ws.Range("H129").Value = 30000
$ws.Range("J129").Value = 30000
$ws.Range("L129").Value = 30000
$ws.Range("N129").Value = -40000

# Row 136
$ws.Range("H136").Value = 864.6
$ws.Range("I136").Value = 775.7143
$ws.Range("J136").Value = 977.7273
$ws.Range("K136").Value = 2327.1429
$ws.Range("L136").Value = 2933.1819
$ws.Range("M136").Value = 222.8571000000002
$ws.Range("N136").Value = -8033.1819
